$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# New header cell N1, matching the format of the existing header row (A1:M1)
$ws.Cells.Item(1, 14).Value = "Correction "
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)

# Blank (but present) data cells N2:N12 under the new column
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 14).Value = "'"
}
